$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.907.22'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +1.47%  '

# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.646.81'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +1.76%  '

# Row 4
$ws.Range('E4').Value = '  -0.06%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '213.56'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.20%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.526'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.46%  '

# Row 7
$ws.Range('E7').Value = '  -0.11%  '

# Row 8
$ws.Range('E8').Value = '  +3.96%  '

# Row 9
$ws.Range('E9').Value = '  +1.19%  '

# Row 10
$ws.Range('E10').Value = '  +0.30%  '

# Row 11
$ws.Range('E11').Value = '  -1.61%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.880.10'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +1.77%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.644.02'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +1.54%  '

# Row 14
$ws.Range('E14').Value = '  +1.52%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.564'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +2.43%  '

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '65.63'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +1.07%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '27.913.09'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +1.41%  '

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '232.19'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +1.06%  '

# Row 19
$ws.Range('E19').Value = '  +1.89%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0₃0723'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.54%  '

# Row 21
$ws.Range('E21').Value = '  -0.03%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '10.73'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +5.61%  '

# Row 23
$ws.Range('E23').Value = '  +2.32%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.16'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +3.32%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '152.19'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.65%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '6.91'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +1.62%  '

# Row 27
$ws.Range('E27').Value = '  +1.17%  '

# Row 28
$ws.Range('E28').Value = '  +0.54%  '

# Row 29
$ws.Range('E29').Value = '  -0.04%  '

# Row 30
$ws.Range('E30').Value = '  +1.72%  '

# Row 31
$ws.Range('E31').Value = '  +0.79%  '

# Row 32
$ws.Range('E32').Value = '  +2.72%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.453.31'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.54%  '

# Row 34
$ws.Range('E34').Value = '  +1.90%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.55'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +1.94%  '

# Row 36
$ws.Range('E36').Value = '  -0.56%  '

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.890'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +3.42%  '

# Row 38
$ws.Range('E38').Value = '  +0.73%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.562'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.10%  '

# Row 40
$ws.Range('E40').Value = '  -1.38%  '

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '69.38'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.35%  '

# Row 42
$ws.Range('E42').Value = '  +2.29%  '

# Row 43
$ws.Range('E43').Value = '  -0.06%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.47'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.52%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.23'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +1.18%  '

# Row 46
$ws.Range('E46').Value = '  -0.49%  '

# Row 47
$ws.Range('E47').Value = '  +5.54%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.788.88'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.65%  '

# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '88.90'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +3.13%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0₆0105'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.34%  '

# Row 51
$ws.Range('E51').Value = '  +1.08%  '
